$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------------
# 1. Capture the existing header-row comments (keyed by their current column
#    number) before we touch anything, since inserting columns later does not
#    automatically relocate cell comments in this runtime.
# ---------------------------------------------------------------------------
$commentCount = $ws.Comments.Count
$commentTexts = @{}
for ($i = 1; $i -le $commentCount; $i++) {
    $cm = $ws.Comments.Item($i)
    $col = $cm.Parent.Column
    $commentTexts[$col] = $cm.Text()
}

# Remove the old comments outright; we will re-create them two columns to the
# right once the new columns are in place.
for ($i = $commentCount; $i -ge 1; $i--) {
    $ws.Comments.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 2. Insert two new blank columns at the very front of the sheet. This shifts
#    the existing header values and every data-validation sqref/formula two
#    columns to the right automatically.
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3. Populate the two new header cells and match the existing header styling
#    (bold, centered, wrapped text - same look as the rest of row 1).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "version"
$ws.Cells.Item(1, 2).Value = "description"

$newHeaders = $ws.Range("A1:B1")
$newHeaders.Font.Bold = $true
$newHeaders.HorizontalAlignment = -4108
$newHeaders.WrapText = $true

# ---------------------------------------------------------------------------
# 4. Re-create the original comments, shifted two columns to the right so
#    they stay attached to the same logical header they originally described.
# ---------------------------------------------------------------------------
foreach ($col in $commentTexts.Keys) {
    $newCol = $col + 2
    $ws.Cells.Item(1, $newCol).AddComment($commentTexts[$col])
}

# New comments describing the two brand-new columns.
$ws.Cells.Item(1, 1).AddComment("Version of the schema to use when validating this metadata.")
$ws.Cells.Item(1, 2).AddComment("Free-text description of this assay.")

# ---------------------------------------------------------------------------
# 5. Add the new "version list" worksheet right after "Export as TSV", with
#    its single allowed value.
# ---------------------------------------------------------------------------
$verSheet = $wb.Worksheets.Add($null, $ws)
$verSheet.Name = "version list"
$verSheet.Range("A1").Value = "'1"

# ---------------------------------------------------------------------------
# 6. Add the data validation for the new "version" column, pointing at the
#    new list sheet.
# ---------------------------------------------------------------------------
$verValidationRange = $ws.Range("A2:A1048576")
$verValidationRange.Validation.Add(3, 1, 1, "'version list'!`$A`$1:`$A`$1")
$verValidationRange.Validation.ErrorTitle = "Value must come from list"
$verValidationRange.Validation.ErrorMessage = "Value must be one of: 1."
$verValidationRange.Validation.IgnoreBlank = $true
$verValidationRange.Validation.InCellDropdown = $true
$verValidationRange.Validation.ShowInput = $true
$verValidationRange.Validation.ShowError = $true
